$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ashley Mueller's affiliations table (Organizations section, row 25) was
# empty; fill in her organization name and its URL.
$ws.Range("A25").Value = "Great Lakes Forestry Centre, Natural Resources Canada"
$ws.Range("B25").Value = "https://www.nrcan.gc.ca/science-data/research-centres-labs/forestry-research-centres/great-lakes-forestry-centre/13459"

# Leave the sheet scrolled/selected on the row that was just edited, like
# the author's saved view (scrolled down with A25:B25 selected).
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("A25:B25").Select()
